$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'59.216.82"
$ws.Range("E2").Formula = "'  +3.18%  "
$ws.Range("D3").Formula = "'2.998.60"
$ws.Range("E3").Formula = "'  +3.34%  "
$ws.Range("E4").Formula = "'  -0.07%  "
$ws.Range("D5").Formula = "'563.91"
$ws.Range("E5").Formula = "'  +2.33%  "
$ws.Range("D6").Formula = "'139.33"
$ws.Range("E6").Formula = "'  +13.32%  "
$ws.Range("E7").Formula = "'  -0.16%  "
$ws.Range("E8").Formula = "'  +4.62%  "
$ws.Range("D9").Formula = "'2.992.05"
$ws.Range("E9").Formula = "'  +3.31%  "
$ws.Range("D10").Formula = "'0.133"
$ws.Range("E10").Formula = "'  +7.72%  "
$ws.Range("E11").Formula = "'  +5.26%  "
$ws.Range("D12").Formula = "'0.458"
$ws.Range("E12").Formula = "'  +4.57%  "
$ws.Range("D13").Formula = "'0.0000231"
$ws.Range("E13").Formula = "'  +9.36%  "
$ws.Range("D14").Formula = "'33.91"
$ws.Range("E14").Formula = "'  +3.64%  "
$ws.Range("E15").Formula = "'  +2.92%  "
$ws.Range("D16").Formula = "'3.493.02"
$ws.Range("E16").Formula = "'  +3.37%  "
$ws.Range("D17").Formula = "'7.03"
$ws.Range("E17").Formula = "'  +7.08%  "
$ws.Range("D18").Formula = "'2.993.88"
$ws.Range("E18").Formula = "'  +3.06%  "
$ws.Range("D19").Formula = "'59.166.16"
$ws.Range("E19").Formula = "'  +3.02%  "
$ws.Range("D20").Formula = "'426.24"
$ws.Range("E20").Formula = "'  +5.89%  "
$ws.Range("D21").Formula = "'13.60"
$ws.Range("E21").Formula = "'  +5.45%  "
$ws.Range("D22").Formula = "'0.714"
$ws.Range("E22").Formula = "'  +6.34%  "
$ws.Range("D23").Formula = "'7.17"
$ws.Range("E23").Formula = "'  +4.94%  "
$ws.Range("D24").Formula = "'13.41"
$ws.Range("E24").Formula = "'  +5.05%  "
$ws.Range("D25").Formula = "'80.55"
$ws.Range("E25").Formula = "'  +4.29%  "
$ws.Range("E26").Formula = "'  -0.02%  "
$ws.Range("E27").Formula = "'  +0.03%  "
$ws.Range("D28").Formula = "'2.14"
$ws.Range("E28").Formula = "'  +11.16%  "
$ws.Range("D29").Formula = "'2.54"
$ws.Range("E29").Formula = "'  +3.53%  "
$ws.Range("D30").Formula = "'7.78"
$ws.Range("E30").Formula = "'  +8.64%  "
$ws.Range("D31").Formula = "'25.66"
$ws.Range("E31").Formula = "'  +3.79%  "
$ws.Range("D32").Formula = "'6.15"
$ws.Range("E32").Formula = "'  +2.21%  "
$ws.Range("D33").Formula = "'0.0993"
$ws.Range("E33").Formula = "'  +0.28%  "
$ws.Range("E34").Formula = "'  +11.39%  "
$ws.Range("D35").Formula = "'0.0" + [char]8323 + "0776"
$ws.Range("E35").Formula = "'  +25.25%  "
$ws.Range("D36").Formula = "'5.77"
$ws.Range("E36").Formula = "'  +6.41%  "
$ws.Range("D37").Formula = "'2.08"
$ws.Range("E37").Formula = "'  +4.42%  "
$ws.Range("D38").Formula = "'49.01"
$ws.Range("E38").Formula = "'  +2.23%  "
$ws.Range("D39").Formula = "'8.68"
$ws.Range("E39").Formula = "'  +5.03%  "
$ws.Range("D40").Formula = "'2.79"
$ws.Range("E40").Formula = "'  +16.14%  "
$ws.Range("D41").Formula = "'406.23"
$ws.Range("E41").Formula = "'  +12.96%  "
$ws.Range("D42").Formula = "'0.0351"
$ws.Range("E42").Formula = "'  +3.42%  "
$ws.Range("D43").Formula = "'2.756.17"
$ws.Range("E43").Formula = "'  +5.11%  "
$ws.Range("E44").Formula = "'  +1.28%  "
$ws.Range("D45").Formula = "'0.247"
$ws.Range("E45").Formula = "'  +7.98%  "
$ws.Range("E46").Formula = "'  +0.02%  "
$ws.Range("D47").Formula = "'125.10"
$ws.Range("E47").Formula = "'  +5.31%  "
$ws.Range("E48").Formula = "'  +4.48%  "
$ws.Range("D49").Formula = "'0.110"
$ws.Range("E49").Formula = "'  +2.72%  "
$ws.Range("D50").Formula = "'32.72"
$ws.Range("E50").Formula = "'  +21.76%  "
$ws.Range("D51").Formula = "'23.54"
$ws.Range("E51").Formula = "'  +3.09%  "
